$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 687.5
$ws.Range("I2").Value = 239.57143
$ws.Range("K2").Value = 239.57143
$ws.Range("M2").Value = -126.57143
$ws.Range("H4").Value = 454.77777
$ws.Range("I4").Value = 206.14285
$ws.Range("J4").Value = 1325
$ws.Range("K4").Value = 206.14285
$ws.Range("L4").Value = 1325
$ws.Range("M4").Value = -92.14285000000001
$ws.Range("N4").Value = -1553
$ws.Range("H9").Value = 33333650
$ws.Range("I9").Value = 66666756
$ws.Range("K9").Value = 66666756
$ws.Range("M9").Value = -66666587
$ws.Range("H17").Value = 2866.6667
$ws.Range("J17").Value = 3750
$ws.Range("L17").Value = 11250
$ws.Range("N17").Value = -11586
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H40").Value = 11830.538
$ws.Range("I40").Value = 5808.4546
$ws.Range("J40").Value = 16246.733
$ws.Range("K40").Value = 5808.4546
$ws.Range("L40").Value = 16246.733
$ws.Range("M40").Value = -5633.4546
$ws.Range("N40").Value = -16596.733
$ws.Range("H113").Value = 4433.0835
$ws.Range("I113").Value = 4133.1113
$ws.Range("K113").Value = 4133.1113
$ws.Range("M113").Value = -879.1112999999996
$ws.Range("H115").Value = 551.75
$ws.Range("I115").Value = 551.75
$ws.Range("K115").Value = 1655.25
$ws.Range("M115").Value = -88.25
$ws.Range("H116").Value = 1596547.8
$ws.Range("I116").Value = 8388.333000000001
$ws.Range("J116").Value = 2787667.2
$ws.Range("K116").Value = 8388.333000000001
$ws.Range("L116").Value = 2787667.2
$ws.Range("M116").Value = -4946.333000000001
$ws.Range("N116").Value = -2794551.2
$ws.Range("H125").Value = 6161.4116
$ws.Range("I125").Value = 4402.5713
$ws.Range("J125").Value = 7392.6
$ws.Range("K125").Value = 39623.14169999999
$ws.Range("L125").Value = 66533.40000000001
$ws.Range("M125").Value = -37163.14169999999
$ws.Range("N125").Value = -71453.40000000001
$ws.Range("H132").Value = 1805.0435
$ws.Range("I132").Value = 1520.6285
$ws.Range("K132").Value = 4561.8855
$ws.Range("M132").Value = -2031.8855
$ws.Range("H137").Value = 338888.06
$ws.Range("I137").Value = 1695.7435
$ws.Range("K137").Value = 5087.2305
$ws.Range("M137").Value = -2537.2305
$ws.Range("H138").Value = 3046
$ws.Range("J138").Value = 3249.0557
$ws.Range("L138").Value = 9747.167099999999
$ws.Range("N138").Value = -20027.1671

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2442.08
$ws.Range("I61").Value = 1322.3684
$ws.Range("K61").Value = 1322.3684
$ws.Range("M61").Value = -1110.3684
$ws.Range("H136").Value = 2442.08
$ws.Range("I136").Value = 1322.3684
$ws.Range("K136").Value = 3967.1052
$ws.Range("M136").Value = -1417.1052

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1931.5625
$ws.Range("I86").Value = 1900.6
$ws.Range("K86").Value = 1900.6
$ws.Range("M86").Value = -777.5999999999999
$ws.Range("H89").Value = 1931.5625
$ws.Range("I89").Value = 1900.6
$ws.Range("K89").Value = 9503
$ws.Range("M89").Value = -3887
$ws.Range("H94").Value = 4309.8667
$ws.Range("I94").Value = 5032.5454
$ws.Range("J94").Value = 2322.5
$ws.Range("K94").Value = 5032.5454
$ws.Range("L94").Value = 2322.5
$ws.Range("M94").Value = -4581.5454
$ws.Range("N94").Value = -3224.5
$ws.Range("H134").Value = 2064.7354
$ws.Range("I134").Value = 1315.8462
$ws.Range("K134").Value = 3947.5386
$ws.Range("M134").Value = -1412.5386

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1525.35
$ws.Range("I16").Value = 1142.5385
$ws.Range("J16").Value = 2236.2856
$ws.Range("K16").Value = 1142.5385
$ws.Range("L16").Value = 2236.2856
$ws.Range("M16").Value = -855.5385000000001
$ws.Range("N16").Value = -2810.2856
$ws.Range("H58").Value = 2536.3333
$ws.Range("I58").Value = 2263
$ws.Range("K58").Value = 2263
$ws.Range("M58").Value = -2060
$ws.Range("H113").Value = 1525.35
$ws.Range("I113").Value = 1142.5385
$ws.Range("J113").Value = 2236.2856
$ws.Range("K113").Value = 1142.5385
$ws.Range("L113").Value = 2236.2856
$ws.Range("M113").Value = 1027.4615
$ws.Range("N113").Value = -6576.2856
$ws.Range("H122").Value = 2282.1892
$ws.Range("I122").Value = 2015.72
$ws.Range("J122").Value = 2837.3333
$ws.Range("K122").Value = 6047.16
$ws.Range("L122").Value = 8511.999899999999
$ws.Range("M122").Value = -3597.16
$ws.Range("N122").Value = -13411.9999
$ws.Range("H136").Value = 2536.3333
$ws.Range("I136").Value = 2263
$ws.Range("K136").Value = 6789
$ws.Range("M136").Value = -4239

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3201.5
$ws.Range("I48").Value = 1403
$ws.Range("K48").Value = 4209
$ws.Range("M48").Value = -3959
$ws.Range("H49").Value = 617.3333
$ws.Range("J49").Value = 537.5
$ws.Range("L49").Value = 1612.5
$ws.Range("N49").Value = -1924.5
$ws.Range("H107").Value = 519.41174
$ws.Range("J107").Value = 324
$ws.Range("L107").Value = 972
$ws.Range("N107").Value = -4812
$ws.Range("H132").Value = 3242.4167
$ws.Range("I132").Value = 2056.5
$ws.Range("K132").Value = 18508.5
$ws.Range("M132").Value = -15978.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13813.25
$ws.Range("I43").Value = 8701.200000000001
$ws.Range("J43").Value = 22333.334
$ws.Range("K43").Value = 8701.200000000001
$ws.Range("L43").Value = 22333.334
$ws.Range("M43").Value = -8550.200000000001
$ws.Range("N43").Value = -22635.334
$ws.Range("H122").Value = 14658.235
$ws.Range("I122").Value = 19957.455
$ws.Range("K122").Value = 59872.36500000001
$ws.Range("M122").Value = -57422.36500000001
$ws.Range("H132").Value = 3831.8333
$ws.Range("I132").Value = 3449.842
$ws.Range("J132").Value = 4147.391
$ws.Range("K132").Value = 10349.526
$ws.Range("L132").Value = 12442.173
$ws.Range("M132").Value = -7819.526
$ws.Range("N132").Value = -17502.173

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7158526.5
$ws.Range("I122").Value = 19416.38
$ws.Range("J122").Value = 28575858
$ws.Range("K122").Value = 58249.14
$ws.Range("L122").Value = 85727574
$ws.Range("M122").Value = -55799.14
$ws.Range("N122").Value = -85732474
$ws.Range("H132").Value = 7430.204
$ws.Range("I132").Value = 9905.969999999999
$ws.Range("J132").Value = 2323.9375
$ws.Range("K132").Value = 29717.91
$ws.Range("L132").Value = 6971.8125
$ws.Range("M132").Value = -27187.91
$ws.Range("N132").Value = -12031.8125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9598.4
$ws.Range("I81").Value = 9623.5
$ws.Range("K81").Value = 19247
$ws.Range("M81").Value = -18186
$ws.Range("H84").Value = 9598.4
$ws.Range("I84").Value = 9623.5
$ws.Range("K84").Value = 96235
$ws.Range("M84").Value = -90931
$ws.Range("H122").Value = 4141.1304
$ws.Range("I122").Value = 3714.5293
$ws.Range("J122").Value = 5349.8335
$ws.Range("K122").Value = 11143.5879
$ws.Range("L122").Value = 16049.5005
$ws.Range("M122").Value = -8693.5879
$ws.Range("N122").Value = -20949.5005
